# "Generate Report for Handback"
# The 0b6255bb-...md file has now been handed back (in sync with en-US),
# so it moves up to the top of the status rows (row 2) and both files show
# "Handed back: in sync with en-US", with refreshed handback timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: rows 2 & 3 swap file identity; both now read
# "Handed back: in sync with en-US" in the zh-cn / de-de columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: rows 2 & 3 swap, status + handback datetime refreshed.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.md"
$wsZhCn.Range("B2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.1292e7eed3d7c840f921ccce1d1471df09ee06a3.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-03-08 02:45:36"
$wsZhCn.Range("E2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.md"
$wsZhCn.Range("F2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.1292e7eed3d7c840f921ccce1d1471df09ee06a3.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-03-08 02:46:19"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.md"
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.f4f9cd88c21682fc64adf2e8b7a1f293b75bd69c.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-08 02:45:36"
$wsZhCn.Range("E3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.md"
$wsZhCn.Range("F3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.f4f9cd88c21682fc64adf2e8b7a1f293b75bd69c.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "2016-03-08 02:46:19"
$wsZhCn.Range("H3").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet: rows 2 & 3 swap, status + handback datetime refreshed.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.md"
$wsDeDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.1292e7eed3d7c840f921ccce1d1471df09ee06a3.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-03-08 02:45:44"
$wsDeDe.Range("E2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.md"
$wsDeDe.Range("F2").Value = "0b6255bb-c6c9-4b2a-8d8f-cfc91798f39f.1292e7eed3d7c840f921ccce1d1471df09ee06a3.de-de.xlf"
$wsDeDe.Range("G2").Value = "2016-03-08 02:46:33"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.md"
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.f4f9cd88c21682fc64adf2e8b7a1f293b75bd69c.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-08 02:45:44"
$wsDeDe.Range("E3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.md"
$wsDeDe.Range("F3").Value = "ea522373-8ce4-4de4-8130-ad3f14132760.f4f9cd88c21682fc64adf2e8b7a1f293b75bd69c.de-de.xlf"
$wsDeDe.Range("G3").Value = "2016-03-08 02:46:33"
$wsDeDe.Range("H3").Value = "Include"
